$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table from row 40 to row 42, copying formatting (style) from the last existing data row
$ws.Range("A40:B40").Copy() | Out-Null
$ws.Range("A41:B42").PasteSpecial(-4122) | Out-Null

# Update product / month grouping rows and their ordered-quantity totals
$ws.Range("A3").Value = "Total"
$ws.Range("B3").Value = 249
$ws.Range("A4").Value = "     [DESK0005] Escritorio personalizable (Personalizado, Blanco)"
$ws.Range("B4").Value = 2
$ws.Range("A5").Value = "          julio 2022"
$ws.Range("B5").Value = 2
$ws.Range("A6").Value = "     [E-COM06] Escritorio de esquina derecho sentarse"
$ws.Range("B6").Value = 10
$ws.Range("A7").Value = "          julio 2022"
$ws.Range("B7").Value = 10
$ws.Range("A8").Value = "     [E-COM07] Gabinete grande"
$ws.Range("B8").Value = 2
$ws.Range("A9").Value = "          julio 2022"
$ws.Range("B9").Value = 1
$ws.Range("A10").Value = "          agosto 2022"
$ws.Range("B10").Value = 1
$ws.Range("A11").Value = "     [E-COM08] Caja de almacenaje"
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = "          julio 2022"
$ws.Range("B12").Value = 10
$ws.Range("A13").Value = "     [E-COM10] Cubo de pedal"
$ws.Range("B13").Value = 25
$ws.Range("A14").Value = "          junio 2022"
$ws.Range("B14").Value = 24
$ws.Range("A15").Value = "          julio 2022"
$ws.Range("B15").Value = 1
$ws.Range("A16").Value = "     [E-COM11] Gabinete con puertas"
$ws.Range("B16").Value = 15
$ws.Range("A17").Value = "          julio 2022"
$ws.Range("B17").Value = 15
$ws.Range("A18").Value = "     [FURN_0096] Escritorio personalizable (Acero, Blanco)"
$ws.Range("B18").Value = 1
$ws.Range("A19").Value = "          julio 2022"
$ws.Range("B19").Value = 1
$ws.Range("A20").Value = "     [FURN_0098] Escritorio personalizable (Aluminio, Blanco)"
$ws.Range("B20").Value = 32
$ws.Range("A21").Value = "          junio 2022"
$ws.Range("B21").Value = 30
$ws.Range("A22").Value = "          agosto 2022"
$ws.Range("B22").Value = 2
$ws.Range("A23").Value = "     [FURN_0269] Silla de oficina negra"
$ws.Range("B23").Value = 27
$ws.Range("A24").Value = "          junio 2022"
$ws.Range("B24").Value = 4
$ws.Range("A25").Value = "          julio 2022"
$ws.Range("B25").Value = 23
$ws.Range("A26").Value = "     [FURN_6666] Pantallas de bloque acústico"
$ws.Range("B26").Value = 83
$ws.Range("A27").Value = "          junio 2022"
$ws.Range("B27").Value = 10
$ws.Range("A28").Value = "          julio 2022"
$ws.Range("B28").Value = 43
$ws.Range("A29").Value = "          agosto 2022"
$ws.Range("B29").Value = 30
$ws.Range("A30").Value = "     [FURN_6741] Gran mesa de reuniones"
$ws.Range("B30").Value = 3
$ws.Range("A31").Value = "          julio 2022"
$ws.Range("B31").Value = 3
$ws.Range("A32").Value = "     [FURN_7777] Silla de oficina"
$ws.Range("B32").Value = 11
$ws.Range("A33").Value = "          junio 2022"
$ws.Range("B33").Value = 4
$ws.Range("A34").Value = "          julio 2022"
$ws.Range("B34").Value = 4
$ws.Range("A35").Value = "          agosto 2022"
$ws.Range("B35").Value = 3
$ws.Range("A36").Value = "     [FURN_8855] Cajón"
$ws.Range("B36").Value = 5
$ws.Range("A37").Value = "          julio 2022"
$ws.Range("B37").Value = 5
$ws.Range("A38").Value = "     [FURN_8888] Lámpara de oficina"
$ws.Range("B38").Value = 22
$ws.Range("A39").Value = "          junio 2022"
$ws.Range("B39").Value = 11
$ws.Range("A40").Value = "          julio 2022"
$ws.Range("B40").Value = 11
$ws.Range("A41").Value = "     [FURN_8999] Sofá de tres asientos"
$ws.Range("B41").Value = 1
$ws.Range("A42").Value = "          julio 2022"
$ws.Range("B42").Value = 1

Write-Output "done"